$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2616.5
$ws.Range("J17").Value = 1998.5555
$ws.Range("L17").Value = 5995.666499999999
$ws.Range("N17").Value = -6331.666499999999

# Row 28
$ws.Range("H28").Value = 586.125
$ws.Range("I28").Value = 312.7143
$ws.Range("J28").Value = 2500
$ws.Range("K28").Value = 312.7143
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 172.2857
$ws.Range("N28").Value = -3470

# Row 33
$ws.Range("H33").Value = 209.84616
$ws.Range("I33").Value = 63.25
$ws.Range("K33").Value = 63.25
$ws.Range("M33").Value = 165.75

# Row 76
$ws.Range("H76").Value = 5923.727
$ws.Range("I76").Value = 4777.25
$ws.Range("K76").Value = 4777.25
$ws.Range("M76").Value = -4462.25

# Row 79
$ws.Range("H79").Value = 5923.727
$ws.Range("I79").Value = 4777.25
$ws.Range("K79").Value = 4777.25
$ws.Range("M79").Value = -3685.25

# Row 132
$ws.Range("H132").Value = 1392.2059
$ws.Range("I132").Value = 1333.9062
$ws.Range("K132").Value = 4001.7186
$ws.Range("M132").Value = -1471.7186

# Row 135
$ws.Range("H135").Value = 883.6957
$ws.Range("I135").Value = 658.381
$ws.Range("K135").Value = 5925.429
$ws.Range("M135").Value = -3390.429

# Row 137
$ws.Range("H137").Value = 1590.9667
$ws.Range("I137").Value = 621.65
$ws.Range("K137").Value = 1864.95
$ws.Range("M137").Value = 685.0500000000002

# Row 138
$ws.Range("H138").Value = 4001.5247
$ws.Range("I138").Value = 1153.75
$ws.Range("J138").Value = 4201.3687
$ws.Range("K138").Value = 3461.25
$ws.Range("L138").Value = 12604.1061
$ws.Range("M138").Value = 1678.75
$ws.Range("N138").Value = -22884.1061

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13666.35
$ws.Range("I32").Value = 5319.675
$ws.Range("J32").Value = 28182.305
$ws.Range("K32").Value = 5319.675
$ws.Range("L32").Value = 28182.305
$ws.Range("M32").Value = -5032.675
$ws.Range("N32").Value = -28756.305

# Row 61
$ws.Range("H61").Value = 1835.1818
$ws.Range("I61").Value = 1837.25
$ws.Range("K61").Value = 1837.25
$ws.Range("M61").Value = -1625.25

# Row 74
$ws.Range("H74").Value = 2952.75
$ws.Range("I74").Value = 974.5833
$ws.Range("K74").Value = 974.5833
$ws.Range("M74").Value = -100.5833

# Row 77
$ws.Range("H77").Value = 2952.75
$ws.Range("I77").Value = 974.5833
$ws.Range("K77").Value = 4872.9165
$ws.Range("M77").Value = -504.9165000000003

# Row 97
$ws.Range("H97").Value = 1333
$ws.Range("I97").Value = 499.5
$ws.Range("K97").Value = 499.5
$ws.Range("M97").Value = -3.5

# Row 122
$ws.Range("H122").Value = 437574.7
$ws.Range("I122").Value = 668284.9399999999
$ws.Range("K122").Value = 2004854.82
$ws.Range("M122").Value = -2002404.82

# Row 136
$ws.Range("H136").Value = 1835.1818
$ws.Range("I136").Value = 1837.25
$ws.Range("K136").Value = 5511.75
$ws.Range("M136").Value = -2961.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2344
$ws.Range("I20").Value = 707.6
$ws.Range("K20").Value = 707.6
$ws.Range("M20").Value = -460.6

# Row 99
$ws.Range("H99").Value = 3775.125
$ws.Range("I99").Value = 3572.8635
$ws.Range("K99").Value = 3572.8635
$ws.Range("M99").Value = -2074.8635

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7526.125
$ws.Range("I31").Value = 3369
$ws.Range("J31").Value = 8120
$ws.Range("K31").Value = 3369
$ws.Range("L31").Value = 8120
$ws.Range("M31").Value = -3074
$ws.Range("N31").Value = -8710

# Row 34
$ws.Range("H34").Value = 7526.125
$ws.Range("I34").Value = 3369
$ws.Range("J34").Value = 8120
$ws.Range("K34").Value = 3369
$ws.Range("L34").Value = 8120
$ws.Range("M34").Value = -3167
$ws.Range("N34").Value = -8524

# Row 94
$ws.Range("H94").Value = 1666.3334
$ws.Range("J94").Value = 1666.3334
$ws.Range("L94").Value = 1666.3334
$ws.Range("N94").Value = -2568.3334

# Row 105
$ws.Range("H105").Value = 1304.2778
$ws.Range("I105").Value = 565.0833
$ws.Range("K105").Value = 565.0833
$ws.Range("M105").Value = 1181.9167

# Row 132
$ws.Range("H132").Value = 2342.5
$ws.Range("I132").Value = 1245.8235
$ws.Range("K132").Value = 3737.4705
$ws.Range("M132").Value = -1207.4705

# Row 134
$ws.Range("H134").Value = 2140.926
$ws.Range("I134").Value = 1219.65
$ws.Range("K134").Value = 3658.95
$ws.Range("M134").Value = -1123.95

# Row 141
$ws.Range("H141").Value = 134938.67
$ws.Range("J141").Value = 134938.67
$ws.Range("L141").Value = 134938.67
$ws.Range("N141").Value = -145298.67

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 195.5
$ws.Range("I8").Value = 195.5
$ws.Range("K8").Value = 586.5
$ws.Range("M8").Value = -447.5

# Row 38
$ws.Range("H38").Value = 68
$ws.Range("I38").Value = 36.875
$ws.Range("K38").Value = 110.625
$ws.Range("M38").Value = 236.375

# Row 117
$ws.Range("H117").Value = 717.7778
$ws.Range("J117").Value = 779.1667
$ws.Range("L117").Value = 2337.5001
$ws.Range("N117").Value = -9221.500100000001

# Row 131
$ws.Range("H131").Value = 2876.2727
$ws.Range("J131").Value = 6703
$ws.Range("L131").Value = 20109
$ws.Range("N131").Value = -30189

# Row 136
$ws.Range("H136").Value = 10900.5
$ws.Range("I136").Value = 7080.6
$ws.Range("K136").Value = 21241.8
$ws.Range("M136").Value = -16141.8

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7326.7144
$ws.Range("I70").Value = 6015
$ws.Range("J70").Value = 8055.4443
$ws.Range("K70").Value = 6015
$ws.Range("L70").Value = 8055.4443
$ws.Range("M70").Value = -5745
$ws.Range("N70").Value = -8595.444299999999

# Row 73
$ws.Range("H73").Value = 7326.7144
$ws.Range("I73").Value = 6015
$ws.Range("J73").Value = 8055.4443
$ws.Range("K73").Value = 6015
$ws.Range("L73").Value = 8055.4443
$ws.Range("M73").Value = -5079
$ws.Range("N73").Value = -9927.444299999999

# Row 97
$ws.Range("H97").Value = 1651.6
$ws.Range("I97").Value = 1544
$ws.Range("J97").Value = 1813
$ws.Range("K97").Value = 1544
$ws.Range("L97").Value = 1813
$ws.Range("M97").Value = -1048
$ws.Range("N97").Value = -2805

# Row 113
$ws.Range("H113").Value = 4360.6665
$ws.Range("I113").Value = 1666.6666
$ws.Range("J113").Value = 5034.1665
$ws.Range("K113").Value = 1666.6666
$ws.Range("L113").Value = 5034.1665
$ws.Range("M113").Value = 503.3334
$ws.Range("N113").Value = -9374.166499999999

# Row 122
$ws.Range("H122").Value = 75432.42999999999
$ws.Range("I122").Value = 4023.1667
$ws.Range("J122").Value = 503888
$ws.Range("K122").Value = 12069.5001
$ws.Range("L122").Value = 1511664
$ws.Range("M122").Value = -9619.500100000001
$ws.Range("N122").Value = -1516564

# Row 132
$ws.Range("H132").Value = 2584.7693
$ws.Range("I132").Value = 1762.6666
$ws.Range("K132").Value = 5287.9998
$ws.Range("M132").Value = -2757.9998

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2497.6667
$ws.Range("I7").Value = 2747
$ws.Range("K7").Value = 2747
$ws.Range("M7").Value = -2635

# Row 16
$ws.Range("H16").Value = 1616.1818
$ws.Range("I16").Value = 1727.8
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 1727.8
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -1557.8
$ws.Range("N16").Value = -840

# Row 126
$ws.Range("H126").Value = 2497.6667
$ws.Range("I126").Value = 2747
$ws.Range("K126").Value = 8241
$ws.Range("M126").Value = -5771

# Row 132
$ws.Range("H132").Value = 3778.7441
$ws.Range("I132").Value = 3431.8076
$ws.Range("K132").Value = 10295.4228
$ws.Range("M132").Value = -7765.4228

# Row 136
$ws.Range("H136").Value = 2046.25
$ws.Range("I136").Value = 2094.7144
$ws.Range("J136").Value = 1707
$ws.Range("K136").Value = 6284.1432
$ws.Range("L136").Value = 5121
$ws.Range("M136").Value = -3734.1432
$ws.Range("N136").Value = -10221

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 781.7273
$ws.Range("I107").Value = 499.8
$ws.Range("J107").Value = 1016.6667
$ws.Range("K107").Value = 1499.4
$ws.Range("L107").Value = 3050.0001
$ws.Range("M107").Value = 420.5999999999999
$ws.Range("N107").Value = -6890.0001

# Row 126
$ws.Range("H126").Value = 3738.6
$ws.Range("I126").Value = 923.25
$ws.Range("K126").Value = 2769.75
$ws.Range("M126").Value = -299.75

# Row 132
$ws.Range("H132").Value = 1259.8966
$ws.Range("I132").Value = 1183.1154
$ws.Range("J132").Value = 1925.3334
$ws.Range("K132").Value = 3549.3462
$ws.Range("L132").Value = 5776.0002
$ws.Range("M132").Value = -1019.3462
$ws.Range("N132").Value = -10836.0002

# Row 136
$ws.Range("H136").Value = 4426.3076
$ws.Range("I136").Value = 1270.5714
$ws.Range("K136").Value = 3811.7142
$ws.Range("M136").Value = -1261.7142
